# Edit script: "inserção da entrada das prestações"
#
# Targets the worked-example block under item (e) "Quantas prestações são
# necessárias para pagar o valor do carro com uma entrada de 3.000?".
# Before:
#   var prestacao = 750
#   var valor = 30000
#   var qtdprest = prestacao / valor
#   qtdprest = 40
# After:
#   var valorTotal = 30000
#   var entrada = 3000
#   var prestacao = 750
#   var qtdPrestfinal = (valorTotal- entrada) / prestacao
#   qtdPrestfinal = 36

$d = $word.ActiveDocument

# Locate the anchor paragraph ("e) Quantas prestações ...") so the four
# following lines can be found unambiguously even if paragraph numbering
# elsewhere in the document shifts.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Quantas*necess*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate anchor paragraph"
}

$p1 = $anchorIndex + 1   # "var prestacao = 750"
$p2 = $anchorIndex + 2   # "var valor = 30000"
$p3 = $anchorIndex + 3   # "var qtdprest = prestacao / valor"
$p4 = $anchorIndex + 4   # "qtdprest = 40"

# 1) "var prestacao = 750" -> "var valorTotal = 30000"
$r1 = $d.Paragraphs.Item($p1).Range
$r1.Find.Execute("prestacao = 750", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "valorTotal = 30000", 2)

# 2) "var valor = 30000" -> "var entrada = 3000"
$r2 = $d.Paragraphs.Item($p2).Range
$r2.Find.Execute("valor = 30000", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "entrada = 3000", 2)

# 3) "var qtdprest = prestacao / valor" -> "var prestacao = 750"
$r3 = $d.Paragraphs.Item($p3).Range
$r3.Find.Execute("qtdprest = prestacao / valor", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "prestacao = 750", 2)

# 4) "qtdprest = 40" -> "var qtdPrestfinal = (valorTotal- entrada) / prestacao"
$r4 = $d.Paragraphs.Item($p4).Range
$r4.Find.Execute("qtdprest = 40", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "var qtdPrestfinal = (valorTotal- entrada) / prestacao", 2)

# 5) Insert a new paragraph after the recalculated line with the final value.
$p4Range = $d.Paragraphs.Item($p4).Range
$p4Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($p4 + 1)
$newPara.Range.Text = "qtdPrestfinal = 36"

Write-Output "done"
